# Add new columns I (I0) and J (IF) to the sheet, populating header + data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: copy formatting from H1 (bold, bordered, centered) to I1:J1 ---
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows: I<row> / J<row> values ---
$data = @(
    @(2, 5, 7),
    @(3, 11, 11),
    @(4, 6, 7),
    @(5, 6, 8),
    @(6, 4, 5),
    @(7, 4, 6),
    @(8, 6, 7),
    @(9, 5, 7),
    @(10, 7, 8),
    @(11, 5, 6),
    @(12, 7, 8),
    @(13, 6, 7),
    @(14, 7, 8),
    @(15, 8, 8),
    @(16, 9, 10),
    @(17, 8, 8),
    @(18, 8, 8),
    @(19, 8, 8),
    @(20, 8, 8),
    @(21, 8, 8),
    @(22, 8, 8),
    @(23, 7, 7),
    @(24, 8, 8),
    @(25, 8, 8),
    @(26, 7, 7),
    @(27, 8, 8),
    @(28, 7, 8),
    @(29, 8, 9),
    @(30, 7, 7),
    @(31, 8, 8),
    @(32, 9, 9),
    @(33, 6, 6),
    @(34, 8, 8),
    @(35, 5, 5),
    @(36, 6, 7),
    @(37, 4, 5),
    @(38, 9, 9),
    @(39, 7, 7),
    @(40, 5, 5),
    @(41, 7, 8),
    @(42, 9, 9),
    @(43, 6, 6),
    @(44, 7, 7)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
